$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Add new row 27: a new task entry for Tran Thanh Tam -----------------
$ws.Range("A27").Value = 15
$ws.Range("B27").Value = "gan lap trinh user admin crud register vo giao dien"
$ws.Range("C27").Value = "Trần Thanh Tâm"
$ws.Range("D27").Value = "'9/6/2023"
$ws.Range("E27").Value = "'9/6/2023"
$ws.Range("F27").Value = 8

# --- 2. Re-create the "Nguyen Minh Tai" block at rows 38-47 (was 30-39) -----
$ws.Range("A38").Value = 1
$ws.Range("B38").Value = "Tìm Kiếm hình ảnh và nội dung,Tài liệu tham khảo"
$ws.Range("C38").Value = "Nguyễn Minh Tài"
$ws.Range("D38").Value = "'3/5/2023"
$ws.Range("E38").Value = "'3/5/2023"
$ws.Range("F38").Value = 2

$ws.Range("A39").Value = 2
$ws.Range("B39").Value = "category Trang chủ "
$ws.Range("C39").Value = "Nguyễn Minh Tài"
$ws.Range("D39").Value = "'4/5/2023"
$ws.Range("E39").Value = "'5/5/2023"
$ws.Range("F39").Value = 2

$ws.Range("A40").Value = 3
$ws.Range("B40").Value = "products "
$ws.Range("C40").Value = "Nguyễn Minh Tài"
$ws.Range("D40").Value = "'7/5/2023"
$ws.Range("E40").Value = "'8/5/2023"
$ws.Range("F40").Value = 2

$ws.Range("A41").Value = 4
$ws.Range("B41").Value = "giỏ hàng"
$ws.Range("C41").Value = "Nguyễn Minh Tài"
$ws.Range("D41").Value = "'10/5/2023"
$ws.Range("E41").Value = "'11/5/2023"
$ws.Range("F41").Value = 2

$ws.Range("A42").Value = 5
$ws.Range("B42").Value = "category hiện thị và xem thông tin sản phẩm"
$ws.Range("C42").Value = "Nguyễn Minh Tài"
$ws.Range("D42").Value = "'12/5/2023"
$ws.Range("E42").Value = "'14/5/2023"
$ws.Range("F42").Value = 3

$ws.Range("A43").Value = 6
$ws.Range("B43").Value = "Trang Chủ"
$ws.Range("C43").Value = "Nguyễn Minh Tài"
$ws.Range("D43").Value = "'15/5/2023"
$ws.Range("E43").Value = "'15/5/2023"
$ws.Range("F43").Value = 2

$ws.Range("A44").Value = 7
$ws.Range("B44").Value = "Profile"
$ws.Range("C44").Value = "Nguyễn Minh Tài"
$ws.Range("D44").Value = "'23/5/2023"
$ws.Range("E44").Value = "'23/5/2023"
$ws.Range("F44").Value = 1

$ws.Range("A45").Value = 8
$ws.Range("B45").Value = " login làm code backend"
$ws.Range("C45").Value = "Nguyễn Minh Tài"
$ws.Range("D45").Value = "'30/5/2023"
$ws.Range("E45").Value = 44963
$ws.Range("F45").Value = 3

$ws.Range("A46").Value = 9
$ws.Range("B46").Value = "lấy dữ liệu và hiển thị thông tin, thêm tạo thông tin"
$ws.Range("C46").Value = "Nguyễn Minh Tài"
$ws.Range("D46").Value = 45083
$ws.Range("E46").Value = 45083
$ws.Range("F46").Value = 7

$ws.Range("A47").Value = 10
$ws.Range("B47").Value = "Thêm nut xóa và nút chỉnh sửa"
$ws.Range("C47").Value = "Nguyễn Minh Tài"
$ws.Range("D47").Value = 45113
$ws.Range("E47").Value = 45113
$ws.Range("F47").Value = 6

# --- 3. Remove the old copy of that block -----------------------------------
# Rows 30-34: clear columns A-F but keep the (already blank) G column cell.
$ws.Range("A30:F34").Clear()

# Rows 35-37: clear entirely so the rows disappear.
$ws.Range("A35:G37").Clear()

# --- 4. Update the view ------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 9
$ws.Range("A38:F47").Select()
